$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph (the 2nd paragraph of the doc),
#    which currently follows the "Play Candy Tower for Free..." H1 heading.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$null = $metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Replace the final paragraph (the italic image-generation prompt) with
#    two paragraphs:
#      - a new bold paragraph: "Play Candy Tower for Free - Exciting Cluster
#        Slot Game"
#      - an (existing-position) italic paragraph whose text becomes the
#        meta-description sentence that used to live near the top of the doc.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$wholeRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Candy Tower for Free - Exciting Cluster Slot Game</w:t></w:r></w:p>' +
            '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Candy Tower, a visually stunning cluster game with a high RTP. Play now for free and activate the exciting Bonus Boost feature.</w:t></w:r></w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$null = $wholeRange.InsertXML($xmlFrag)
